# Edit script: renames the "Requested quantity" headers on the two
# existing sheets and adds a new "PO Forecast" sheet (with ds,
# PO_Forecast, yhat_lower, yhat_upper columns) after them.

$wb = $excel.ActiveWorkbook
$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- 1. Add the new "PO Forecast" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "PO Forecast"

# Match the sheetPr/pageMargins conventions used by the other sheets.
$ws3.Outline.SummaryRow = 1
$ws3.Outline.SummaryColumn = 1
$ws3.PageSetup.LeftMargin = 54
$ws3.PageSetup.RightMargin = 54
$ws3.PageSetup.TopMargin = 72
$ws3.PageSetup.BottomMargin = 72
$ws3.PageSetup.HeaderMargin = 36
$ws3.PageSetup.FooterMargin = 36

# --- 2. Header row: reuse the bold/bordered header style from "Weekly Quantity" ---
$wsWeekly.Range("A1:B1").Copy($ws3.Range("A1:B1"))
$wsWeekly.Range("A1:B1").Copy($ws3.Range("C1:D1"))
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# --- 3. Data rows: reuse the date-formatted style from column A on "Weekly Quantity" ---
for ($r = 2; $r -le 38; $r++) {
    $wsWeekly.Range("A2").Copy($ws3.Cells.Item($r, 1))
}

# --- 4. Fill in the forecast values ---
$ws3.Cells.Item(2, 1).Value = 45032.99999999999
$ws3.Cells.Item(2, 2).Value = 33
$ws3.Cells.Item(2, 3).Value = -149.4623435216175
$ws3.Cells.Item(2, 4).Value = 220.4587268990612
$ws3.Cells.Item(3, 1).Value = 45039.99999999999
$ws3.Cells.Item(3, 2).Value = 39
$ws3.Cells.Item(3, 3).Value = -158.6019897442179
$ws3.Cells.Item(3, 4).Value = 214.6034755328131
$ws3.Cells.Item(4, 1).Value = 45046.99999999999
$ws3.Cells.Item(4, 2).Value = 45
$ws3.Cells.Item(4, 3).Value = -148.0569877925644
$ws3.Cells.Item(4, 4).Value = 234.1487066376461
$ws3.Cells.Item(5, 1).Value = 45060.99999999999
$ws3.Cells.Item(5, 2).Value = 56
$ws3.Cells.Item(5, 3).Value = -138.8115313470916
$ws3.Cells.Item(5, 4).Value = 230.4975436591615
$ws3.Cells.Item(6, 1).Value = 45081.99999999999
$ws3.Cells.Item(6, 2).Value = 74
$ws3.Cells.Item(6, 3).Value = -112.0194308763698
$ws3.Cells.Item(6, 4).Value = 259.6443069342743
$ws3.Cells.Item(7, 1).Value = 45088.99999999999
$ws3.Cells.Item(7, 2).Value = 79
$ws3.Cells.Item(7, 3).Value = -106.6238841509041
$ws3.Cells.Item(7, 4).Value = 265.356714584919
$ws3.Cells.Item(8, 1).Value = 45095.99999999999
$ws3.Cells.Item(8, 2).Value = 85
$ws3.Cells.Item(8, 3).Value = -118.1467638216152
$ws3.Cells.Item(8, 4).Value = 277.9227004925371
$ws3.Cells.Item(9, 1).Value = 45102.99999999999
$ws3.Cells.Item(9, 2).Value = 91
$ws3.Cells.Item(9, 3).Value = -106.4994379431113
$ws3.Cells.Item(9, 4).Value = 279.3870107131009
$ws3.Cells.Item(10, 1).Value = 45109.99999999999
$ws3.Cells.Item(10, 2).Value = 97
$ws3.Cells.Item(10, 3).Value = -93.79326059505452
$ws3.Cells.Item(10, 4).Value = 285.4986744689729
$ws3.Cells.Item(11, 1).Value = 45123.99999999999
$ws3.Cells.Item(11, 2).Value = 108
$ws3.Cells.Item(11, 3).Value = -75.75229373030483
$ws3.Cells.Item(11, 4).Value = 303.2077157124473
$ws3.Cells.Item(12, 1).Value = 45130.99999999999
$ws3.Cells.Item(12, 2).Value = 114
$ws3.Cells.Item(12, 3).Value = -77.79880156914601
$ws3.Cells.Item(12, 4).Value = 296.3390171572648
$ws3.Cells.Item(13, 1).Value = 45144.99999999999
$ws3.Cells.Item(13, 2).Value = 125
$ws3.Cells.Item(13, 3).Value = -68.41985334433502
$ws3.Cells.Item(13, 4).Value = 318.8413199249561
$ws3.Cells.Item(14, 1).Value = 45151.99999999999
$ws3.Cells.Item(14, 2).Value = 131
$ws3.Cells.Item(14, 3).Value = -71.39804344541351
$ws3.Cells.Item(14, 4).Value = 313.115268442962
$ws3.Cells.Item(15, 1).Value = 45172.99999999999
$ws3.Cells.Item(15, 2).Value = 148
$ws3.Cells.Item(15, 3).Value = -56.36797429134486
$ws3.Cells.Item(15, 4).Value = 333.8993774498993
$ws3.Cells.Item(16, 1).Value = 45179.99999999999
$ws3.Cells.Item(16, 2).Value = 154
$ws3.Cells.Item(16, 3).Value = -23.26703110427812
$ws3.Cells.Item(16, 4).Value = 338.6862903913225
$ws3.Cells.Item(17, 1).Value = 45186.99999999999
$ws3.Cells.Item(17, 2).Value = 160
$ws3.Cells.Item(17, 3).Value = -25.23985515466447
$ws3.Cells.Item(17, 4).Value = 337.5803137425985
$ws3.Cells.Item(18, 1).Value = 45193.99999999999
$ws3.Cells.Item(18, 2).Value = 165
$ws3.Cells.Item(18, 3).Value = -24.92898777471201
$ws3.Cells.Item(18, 4).Value = 354.3614239416337
$ws3.Cells.Item(19, 1).Value = 45200.99999999999
$ws3.Cells.Item(19, 2).Value = 171
$ws3.Cells.Item(19, 3).Value = -18.8172178471034
$ws3.Cells.Item(19, 4).Value = 354.3919218302834
$ws3.Cells.Item(20, 1).Value = 45207.99999999999
$ws3.Cells.Item(20, 2).Value = 177
$ws3.Cells.Item(20, 3).Value = -5.564315856688926
$ws3.Cells.Item(20, 4).Value = 360.3327699807875
$ws3.Cells.Item(21, 1).Value = 45214.99999999999
$ws3.Cells.Item(21, 2).Value = 183
$ws3.Cells.Item(21, 3).Value = -3.690611064769317
$ws3.Cells.Item(21, 4).Value = 373.1369261810859
$ws3.Cells.Item(22, 1).Value = 45221.99999999999
$ws3.Cells.Item(22, 2).Value = 188
$ws3.Cells.Item(22, 3).Value = -6.696714850926763
$ws3.Cells.Item(22, 4).Value = 367.6618083860693
$ws3.Cells.Item(23, 1).Value = 45228.99999999999
$ws3.Cells.Item(23, 2).Value = 194
$ws3.Cells.Item(23, 3).Value = 1.131033061522316
$ws3.Cells.Item(23, 4).Value = 381.1961279625451
$ws3.Cells.Item(24, 1).Value = 45235.99999999999
$ws3.Cells.Item(24, 2).Value = 200
$ws3.Cells.Item(24, 3).Value = 20.27579993027134
$ws3.Cells.Item(24, 4).Value = 381.4796064245278
$ws3.Cells.Item(25, 1).Value = 45249.99999999999
$ws3.Cells.Item(25, 2).Value = 211
$ws3.Cells.Item(25, 3).Value = 20.88885933142781
$ws3.Cells.Item(25, 4).Value = 391.2616662914564
$ws3.Cells.Item(26, 1).Value = 45256.99999999999
$ws3.Cells.Item(26, 2).Value = 217
$ws3.Cells.Item(26, 3).Value = 28.74245115048847
$ws3.Cells.Item(26, 4).Value = 414.912676097954
$ws3.Cells.Item(27, 1).Value = 45270.99999999999
$ws3.Cells.Item(27, 2).Value = 228
$ws3.Cells.Item(27, 3).Value = 21.51087407188218
$ws3.Cells.Item(27, 4).Value = 410.7345002295195
$ws3.Cells.Item(28, 1).Value = 45277.99999999999
$ws3.Cells.Item(28, 2).Value = 234
$ws3.Cells.Item(28, 3).Value = 41.63918609011267
$ws3.Cells.Item(28, 4).Value = 424.1737927287734
$ws3.Cells.Item(29, 1).Value = 45298.99999999999
$ws3.Cells.Item(29, 2).Value = 251
$ws3.Cells.Item(29, 3).Value = 68.25752869139443
$ws3.Cells.Item(29, 4).Value = 466.3342125287227
$ws3.Cells.Item(30, 1).Value = 45312.99999999999
$ws3.Cells.Item(30, 2).Value = 263
$ws3.Cells.Item(30, 3).Value = 73.08378488587995
$ws3.Cells.Item(30, 4).Value = 454.8103464090547
$ws3.Cells.Item(31, 1).Value = 45319.99999999999
$ws3.Cells.Item(31, 2).Value = 269
$ws3.Cells.Item(31, 3).Value = 83.35628242987119
$ws3.Cells.Item(31, 4).Value = 455.3601507993885
$ws3.Cells.Item(32, 1).Value = 45326.99999999999
$ws3.Cells.Item(32, 2).Value = 274
$ws3.Cells.Item(32, 3).Value = 91.81474286279384
$ws3.Cells.Item(32, 4).Value = 461.8261256355665
$ws3.Cells.Item(33, 1).Value = 45333.99999999999
$ws3.Cells.Item(33, 2).Value = 280
$ws3.Cells.Item(33, 3).Value = 106.2094547635467
$ws3.Cells.Item(33, 4).Value = 463.5309336553689
$ws3.Cells.Item(34, 1).Value = 45340.99999999999
$ws3.Cells.Item(34, 2).Value = 286
$ws3.Cells.Item(34, 3).Value = 89.47483637219852
$ws3.Cells.Item(34, 4).Value = 465.5220779652544
$ws3.Cells.Item(35, 1).Value = 45347.99999999999
$ws3.Cells.Item(35, 2).Value = 292
$ws3.Cells.Item(35, 3).Value = 93.56535860533052
$ws3.Cells.Item(35, 4).Value = 482.9278561931833
$ws3.Cells.Item(36, 1).Value = 45354.99999999999
$ws3.Cells.Item(36, 2).Value = 297
$ws3.Cells.Item(36, 3).Value = 110.0504894882028
$ws3.Cells.Item(36, 4).Value = 485.331372175574
$ws3.Cells.Item(37, 1).Value = 45361.99999999999
$ws3.Cells.Item(37, 2).Value = 303
$ws3.Cells.Item(37, 3).Value = 122.6032585174171
$ws3.Cells.Item(37, 4).Value = 490.6465701566764
$ws3.Cells.Item(38, 1).Value = 45368.99999999999
$ws3.Cells.Item(38, 2).Value = 309
$ws3.Cells.Item(38, 3).Value = 114.8162591996015
$ws3.Cells.Item(38, 4).Value = 502.5113371826332

# --- 5. Rename the "Requested quantity" headers on the existing sheets ---
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

$ws3.Range("A1").Select()
